# Update gh-pages output data for 合肥-漫展信息 workbook.
# Two sheets ("展览" and "全部类型") each contain a row for the same events;
# the "想去人数" (wanted-to-go count) column F needs to be refreshed.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 5549
    $ws.Range("F7").Value = 42

    if ($name -eq "展览") {
        $ws.Range("F8").Value = 365
    } else {
        $ws.Range("F9").Value = 365
    }
}
